$p = $ppt.ActivePresentation

# Slide 11 ("The Dataset") - content placeholder, first bullet:
# "It contains # columns" -> "It contains # columns222"
$slide = $p.Slides.Item(11)
$shape = $slide.Shapes.Item(2)
$textRange = $shape.TextFrame.TextRange
$paragraph = $textRange.Paragraphs(1, 1)

# Replace just the "# columns" portion so the run splits the same way a
# user editing in-place would split it: "It contains " + "# columns222"
$target = $paragraph.Characters(13, 9)
$target.Text = "# columns222"
